$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone row 11 (values + formatting) down into the new row 12, then
# overwrite the cells that actually differ. Copying first keeps the
# per-column cell styles (s="1" for A/F, s="2" for C/D, default for B/E)
# identical to the rest of the table instead of Excel minting new xf
# records the way a bare `.Style =` assignment would.
$ws.Range("A11:F11").Copy($ws.Range("A12:F12"))

$ws.Range("A12").Value = "com.singleton.strechy"
$ws.Range("B12").Value = "stretchy"
$ws.Range("C12").Value = "shmualtamara@gmail.com"
$ws.Range("D12").Value = "shmulmaor2@gmail.com"
$ws.Range("E12").Value = "27/5/2019 15:59"
$ws.Range("F12").Value = "Fantastic app with great thinking behind it. Addictive as hell"

# Add the mailto hyperlink for the recovery-email column only (matches
# the source row layout where D has a live hyperlink and C does not).
$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:shmulmaor2@gmail.com", [Type]::Missing, [Type]::Missing, "shmulmaor2@gmail.com") | Out-Null

# Hyperlinks.Add re-stamps the cell with Excel's built-in "Hyperlink"
# style; restore the plain centered style shared by the rest of column D
# (and re-apply the value, since Copy also carries the source value).
$ws.Range("D11").Copy($ws.Range("D12"))
$ws.Range("D12").Value = "shmulmaor2@gmail.com"

$ws.Range("F12").Select()
